$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 10-13 (Resolving-Mac as Sending cluster) - data no longer present after TPM update
$ws.Rows("10:13").Delete()

# Update numeric values (columns E:T) for rows 2-9 with refreshed TPM-derived statistics
$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$row2 = @(3, 1, 23.544642, 70.633926, 0.8903342714957673, 0.8903342714957674, 3, 1, 24.576554, 73.729662, 0.07553767049546639, 0.07553767049546638, 578.6461655236681, 5207.815489713013, 0.06725377683106838, 0.06725377683106838)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "2").Value = $row2[$i] }

$row3 = @(3, 1, 23.544642, 70.633926, 0.8903342714957673, 0.8903342714957674, 3, 1, 117.0512696666667, 351.153809, 0.359764849016532, 0.359764849016532, 2755.930239947126, 24803.37215952414, 0.3203109747589187, 0.3203109747589187)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "3").Value = $row3[$i] }

$row4 = @(3, 1, 23.544642, 70.633926, 0.8903342714957673, 0.8903342714957674, 3, 1, 55.68784966666667, 167.063549, 0.1711603033819035, 0.1711603033819035, 1311.150484151486, 11800.35435736338, 0.1523898840205216, 0.1523898840205216)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "4").Value = $row4[$i] }

$row5 = @(3, 1, 23.544642, 70.633926, 0.8903342714957673, 0.8903342714957674, 3, 1, 128.0392633333333, 384.11779, 0.3935371771060981, 0.3935371771060981, 3014.63861712706, 27131.74755414354, 0.3503796358852586, 0.3503796358852586)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "5").Value = $row5[$i] }

$row6 = @(3, 1, 2.900079666666667, 8.700239, 0.1096657285042327, 0.1096657285042327, 3, 1, 24.576554, 73.729662, 0.07553767049546639, 0.07553767049546638, 71.27396453213534, 641.465680789218, 0.008283893664398006, 0.008283893664398004)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "6").Value = $row6[$i] }

$row7 = @(3, 1, 2.900079666666667, 8.700239, 0.1096657285042327, 0.1096657285042327, 3, 1, 117.0512696666667, 351.153809, 0.359764849016532, 0.359764849016532, 339.4580071178168, 3055.122064060351, 0.03945387425761326, 0.03945387425761326)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "7").Value = $row7[$i] }

$row8 = @(3, 1, 2.900079666666667, 8.700239, 0.1096657285042327, 0.1096657285042327, 3, 1, 55.68784966666667, 167.063549, 0.1711603033819035, 0.1711603033819035, 161.4992004986901, 1453.492804488211, 0.01877041936138193, 0.01877041936138193)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "8").Value = $row8[$i] }

$row9 = @(3, 1, 2.900079666666667, 8.700239, 0.1096657285042327, 0.1096657285042327, 3, 1, 128.0392633333333, 384.11779, 0.3935371771060981, 0.3935371771060981, 371.3240641279789, 3341.91657715181, 0.04315754122083949, 0.04315754122083949)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "9").Value = $row9[$i] }
